$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: the "Tonga Energy Road Map" sentence is reassessed from "yes" (relevant,
#     with mitigation/emissions/policy-target coding) to "no" (not relevant) - the
#     detailed coding columns are cleared. ---
$ws.Range("C12:H12").ClearContents()
$ws.Range("B12").Value = "no"

# --- Row 16: the "Global Commons" sentence keeps its "yes" coding, but the Unit is
#     recoded from "measures" to "n.a.", and the 30-word explanation is reworded. ---
$ws.Range("D16").Value = "n.a."
$ws.Range("H16").Value = "the co-existence in a shared world where we should work together to improve this world for the benefit of all "

# --- Row 17: the "right to exist" sentence is now reassessed from "no" to "yes",
#     with a full set of coding columns added. ---
$ws.Range("B17").Value = "yes"
$ws.Range("C17").Value = "other(right to exist)"
$ws.Range("D17").Value = "n.a."
$ws.Range("E17").Value = "national"
$ws.Range("F17").Value = "n.a."
$ws.Range("G17").Value = "sufficientarian"
$ws.Range("H17").Value = "moral judgement on the right to exist. Setting a minimum threshold. "

# --- Row 18: Unit recoded from "measures" to "n.a.". ---
$ws.Range("D18").Value = "n.a."

# --- Row 20: Topic updated to add "implementation". ---
$ws.Range("C20").Value = "moral responsiblity, implementation"

# Move the active selection to reflect where the edits were focused.
$ws.Range("C16").Select()
